$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"

# Column D holds a date value, formatted like the row above it (style index 2 / YYYY-MM-DD HH:MM:SS)
# Set the number format first, then assign the raw date serial (44656 => 2022-04-05)
# so Excel doesn't manufacture a brand-new/unused number format style.
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(56, 4).NumberFormat
$ws.Cells.Item($row, 4).Value = 44656

$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 100112031
$ws.Cells.Item($row, 7).Value = "Poroto verde"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 1500
$ws.Cells.Item($row, 11).Value = 600
$ws.Cells.Item($row, 12).Value = 700
$ws.Cells.Item($row, 13).Value = 650
$ws.Cells.Item($row, 14).Value = "$/kilo"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 650
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"

$wb.Save()
